$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values
$ws.Range("B2").Value = 516.92082400000004
$ws.Range("C2").Value = 435.435768
$ws.Range("D2").Value = 516.92082400000004
$ws.Range("E2").Value = 435.435768

# Row 3 values - C3 removed, D3 added, B3/E3 updated
$ws.Range("B3").Value = 524.56004800000005
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 524.56004800000005
$ws.Range("E3").Value = 435.435768

# Selection change
$ws.Range("B1:E3").Select()
